# Applies the cryptos.xlsx price/volume refresh described in the commit
# "Updated cryptos list on Tue Dec 19 23:43:50 UTC 2023 with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.395.13'
$ws.Range('E2').Value = '  -0.70%  '
$ws.Range('D3').Value = '2.179.35'
$ws.Range('E3').Value = '  -1.70%  '
$ws.Range('E4').Value = '  -0.17%  '
$ws.Range('D5').Formula = "'252.70"
$ws.Range('E5').Value = '  +4.90%  '
$ws.Range('E6').Value = '  -1.48%  '
$ws.Range('D7').Formula = "'73.10"
$ws.Range('E7').Value = '  -2.27%  '
$ws.Range('E8').Value = '  -0.08%  '
$ws.Range('D9').Formula = "'0.582"
$ws.Range('E9').Value = '  -3.90%  '
$ws.Range('D10').Formula = "'40.08"
$ws.Range('E10').Value = '  -2.68%  '
$ws.Range('D11').Formula = "'0.0903"
$ws.Range('E11').Value = '  -2.83%  '
$ws.Range('E12').Value = '  -0.33%  '
$ws.Range('D13').Formula = "'6.73"
$ws.Range('E13').Value = '  -1.99%  '
$ws.Range('D14').Value = '2.506.03'
$ws.Range('E14').Value = '  -1.60%  '
$ws.Range('D15').Formula = "'14.18"
$ws.Range('E15').Value = '  -3.65%  '
$ws.Range('D16').Value = '2.172.53'
$ws.Range('E16').Value = '  -1.92%  '
$ws.Range('D17').Formula = "'0.764"
$ws.Range('E17').Value = '  -4.76%  '
$ws.Range('D18').Value = '42.295.11'
$ws.Range('E18').Value = '  -0.60%  '
$ws.Range('E19').Value = '  -2.79%  '
$ws.Range('D20').Formula = "'70.53"
$ws.Range('E20').Value = '  -0.35%  '
$ws.Range('D21').Formula = "'5.83"
$ws.Range('E21').Value = '  -1.58%  '
$ws.Range('D22').Formula = "'226.75"
$ws.Range('E22').Value = '  -1.25%  '
$ws.Range('D23').Formula = "'9.36"
$ws.Range('E24').Value = '  -2.51%  '
$ws.Range('E25').Value = '  +0.06%  '
$ws.Range('D26').Formula = "'10.44"
$ws.Range('E26').Value = '  -4.42%  '
$ws.Range('D27').Formula = "'3.37"
$ws.Range('E27').Value = '  -0.12%  '
$ws.Range('E28').Value = '  -2.41%  '
$ws.Range('E29').Value = '  +2.01%  '
$ws.Range('D30').Formula = "'36.77"
$ws.Range('E30').Value = '  -0.23%  '
$ws.Range('D31').Formula = "'170.10"
$ws.Range('E31').Value = '  -1.62%  '
$ws.Range('E32').Value = '  -1.57%  '
$ws.Range('D33').Formula = "'0.0814"
$ws.Range('E33').Value = '  +2.25%  '
$ws.Range('E34').Value = '  -3.55%  '
$ws.Range('E35').Value = '  -1.64%  '
$ws.Range('D36').Formula = "'0.107"
$ws.Range('E36').Value = '  -1.82%  '
$ws.Range('E37').Value = '  +3.56%  '
$ws.Range('D38').Formula = "'4.19"
$ws.Range('E38').Value = '  -5.04%  '
$ws.Range('D39').Formula = "'11.74"
$ws.Range('E39').Value = '  -5.19%  '
$ws.Range('E40').Value = '  -4.04%  '
$ws.Range('B41').Value = 'Algorand'
$ws.Range('C41').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D41').Formula = "'0.195"
$ws.Range('E41').Value = '  -1.20%  '
$ws.Range('B42').Value = 'MultiversX'
$ws.Range('C42').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D42').Formula = "'59.22"
$ws.Range('E42').Value = '  -1.95%  '
$ws.Range('D43').Formula = "'5.14"
$ws.Range('E43').Value = '  -6.68%  '
$ws.Range('D44').Formula = "'101.48"
$ws.Range('E44').Value = '  +2.22%  '
$ws.Range('D45').Formula = "'2.45"
$ws.Range('E45').Value = '  +7.52%  '
$ws.Range('E46').Value = '  -1.87%  '
$ws.Range('D47').Formula = "'0.458"
$ws.Range('E47').Value = '  +6.97%  '
$ws.Range('E48').Value = '  -5.23%  '
$ws.Range('E49').Value = '  -2.36%  '
$ws.Range('E50').Value = '  -1.45%  '
$ws.Range('E51').Value = '  +0.25%  '
